# Update "想去人数" (want-to-go count) values in 展览 (sheet1) and 全部类型 (sheet4) sheets
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws4 = $wb.Worksheets.Item(4)

$ws1.Range("F3").Value2 = 586
$ws1.Range("F5").Value2 = 299
$ws1.Range("F7").Value2 = 1459
$ws1.Range("F10").Value2 = 758
$ws1.Range("F12").Value2 = 180
$ws1.Range("F14").Value2 = 457
$ws1.Range("F15").Value2 = 1401
$ws1.Range("F16").Value2 = 130
$ws1.Range("F18").Value2 = 283
$ws1.Range("F20").Value2 = 77
$ws1.Range("F21").Value2 = 664
$ws1.Range("F22").Value2 = 1017
$ws1.Range("F23").Value2 = 40
$ws1.Range("F24").Value2 = 251
$ws1.Range("F25").Value2 = 29
$ws1.Range("F26").Value2 = 6032
$ws1.Range("F27").Value2 = 73
$ws1.Range("F29").Value2 = 119
$ws1.Range("F31").Value2 = 14740
$ws1.Range("F32").Value2 = 1465
$ws1.Range("F33").Value2 = 236
$ws1.Range("F34").Value2 = 108
$ws1.Range("F36").Value2 = 10118
$ws1.Range("F37").Value2 = 660
$ws1.Range("F38").Value2 = 4231
$ws1.Range("F39").Value2 = 166
$ws1.Range("F41").Value2 = 113

$ws4.Range("F3").Value2 = 586
$ws4.Range("F5").Value2 = 299
$ws4.Range("F7").Value2 = 1459
$ws4.Range("F10").Value2 = 758
$ws4.Range("F12").Value2 = 180
$ws4.Range("F14").Value2 = 457
$ws4.Range("F15").Value2 = 1401
$ws4.Range("F16").Value2 = 130
$ws4.Range("F18").Value2 = 283
$ws4.Range("F21").Value2 = 77
$ws4.Range("F22").Value2 = 664
$ws4.Range("F24").Value2 = 1017
$ws4.Range("F26").Value2 = 251
$ws4.Range("F27").Value2 = 29
$ws4.Range("F29").Value2 = 6032
$ws4.Range("F30").Value2 = 73
$ws4.Range("F32").Value2 = 119
$ws4.Range("F34").Value2 = 14740
$ws4.Range("F35").Value2 = 1465
$ws4.Range("F36").Value2 = 236
$ws4.Range("F37").Value2 = 108
$ws4.Range("F39").Value2 = 10118
$ws4.Range("F40").Value2 = 660
$ws4.Range("F41").Value2 = 4231
$ws4.Range("F42").Value2 = 166
$ws4.Range("F44").Value2 = 113
